# Apply updated cryptos list values (price/volume refresh + two coin swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.232.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.63%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.135.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.72%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.133.07'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.22%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.75%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.479'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.78%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000260'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.68%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.602.64'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.99%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.308.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.65%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.114'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.125.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.18%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.31%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '483.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.54%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.713'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.78'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.92%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.63%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.86'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.86%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.66%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.25%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.125'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.72'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -7.51%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.85'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.58%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.92%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.68%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +7.51%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '54.68'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.32%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0749'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '450.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -10.46%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.126'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.90%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0403'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.29%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.10%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.887.53'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.81%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.275'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.47%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.99'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.116'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.51%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.50'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.45%  '
$ws.Range('E51').Style = 'Normal'
